$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '24.839.24'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +0.80%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.703.90'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +0.25%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '316.47'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.31%  '

$ws.Range("E6").Value = '  +0.11%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3935'
$ws.Range("D7").ClearFormats()

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4049'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +0.20%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.518'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -1.71%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.002'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +0.09%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '53.61'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -2.66%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08893'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +0.91%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.338'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +0.81%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '23.56'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +0.52%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.036'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +5.13%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001332'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -0.03%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.703.11'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +0.11%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '100.32'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -0.97%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.07048'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -0.59%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '19.72'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -0.22%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.073'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +2.03%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '14.61'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +3.04%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '24.813.48'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.76%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.199'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +7.37%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.362'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +0.97%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.86'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +1.87%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '162.15'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +1.53%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.433'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +10.73%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '136.72'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +2.03%  '

$ws.Range("E31").Value = '  -1.71%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.08912'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +3.89%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '11.20'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -0.98%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.980'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +1.22%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.2754'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -0.51%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '14.48'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -2.00%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.09234'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +1.69%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.02764'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -2.07%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.463'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -0.14%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.7728'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -0.44%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '15.93'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +1.77%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.7219'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -0.95%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.578'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +2.31%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.210'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -0.33%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.001'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +0.09%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '140.71'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -0.82%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.325'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -4.79%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '91.15'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +3.15%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.08004'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -0.45%  '

$ws.Range("B33").Value = 'ImmutableX'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.086'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -1.76%  '

$ws.Range("B34").Value = 'InternetComputer(DFINITY)'
$ws.Range("C34").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '7.473'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +0.98%  '
